# Stage 3 Analysis Complete
# Update the "19812000" sheet (stage 2 -> stage 2 + stage 3 dominant-frequency
# analysis): retarget the intensity headers to median/IQR, add dominant
# frequency median/IQR/AD-test/n columns to the summary table, refresh the
# recomputed intensity statistics, and append the new one-way ANOVA +
# Tukey-Kramer tables for dominant frequency.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("19812000")
$ws.Activate()

# --- Row 1: header row (A1:L1) ------------------------------------------
# E1/F1 change from "average/std" to "median/IQR" wording; I1/J1 likewise
# become "Dominant Frequency, median/IQR"; K1/L1 are brand new columns
# (AD test, normality / n) mirroring D1/H1.
$ws.Range("E1").Value = "Intensity (mV^2/s), median"
$ws.Range("F1").Value = "Intensity (mV^2/s), IQR"
$ws.Range("I1").Value = "Dominant Frequency, median"
$ws.Range("J1").Value = "Dominant Frequency, IQR"
$ws.Range("K1").Value = "AD test, normality"
$ws.Range("L1").Value = "n"

# --- Row 2: Control ---------------------------------------------------
$ws.Range("E2").Value = 378566.68421590159
$ws.Range("F2").Value = 346753.61093470646
$ws.Range("G2").Value = 0.016455119890529613
$ws.Range("I2").Value = 15.418228149414063
$ws.Range("J2").Value = 16.054031372070313
$ws.Range("K2").Value = 0.0005
$ws.Range("L2").Value = 28

# --- Row 3: Control, pH 6.9 --------------------------------------------
$ws.Range("E3").Value = 414085.71964627929
$ws.Range("F3").Value = 125928.82214807905
$ws.Range("G3").Value = 0.57979270012249917
$ws.Range("I3").Value = 16.212982177734375
$ws.Range("J3").Value = 7.8680648803710938
$ws.Range("K3").Value = 0.065940482305243589
$ws.Range("L3").Value = 33

# --- Row 4: treatment group renamed "Washout" -> "PostTest", plus refreshed
#     stats and the new dominant-frequency columns -----------------------
$ws.Range("A4").Value = "PostTest"
$ws.Range("E4").Value = 295410.50212958449
$ws.Range("F4").Value = 72611.518940786162
$ws.Range("G4").Value = 0.38094222869523064
$ws.Range("I4").Value = 18.120391845703125
$ws.Range("J4").Value = 6.1196060180664063
$ws.Range("K4").Value = 0.97235305289678287
$ws.Range("L4").Value = 11

# --- Row 5: KS Test / Cliffs D sub-header, mirrored into I5:K5 ----------
$ws.Range("I5").Value = "p-value"
$ws.Range("J5").Value = "KS D stat"
$ws.Range("K5").Value = "Cliffs D"

# --- Row 6: KS Test, 1 vs 2 stats, plus the new dominant-frequency values
$ws.Range("G6").Value = -0.082251082251082255
$ws.Range("I6").Value = 0.068571263419851772
$ws.Range("J6").Value = 0.3214285714285714
$ws.Range("K6").Value = 0.21536796536796537

# --- one-way ANOVA, intensity (rows 20-24) - recomputed numbers --------
$ws.Range("B22").Value = 120001217405.62329
$ws.Range("D22").Value = 60000608702.811646
$ws.Range("E22").Value = 3.8279595183402142
$ws.Range("F22").Value = 0.026513505941062845

$ws.Range("B23").Value = 1081527111417.6011
$ws.Range("D23").Value = 15674305962.573929

$ws.Range("B24").Value = 1201528328823.2244

# --- Multiple Comparison (Tukey-Kramer method), intensity (rows 26-30) -
$ws.Range("C28").Value = -124268.14471554037
$ws.Range("D28").Value = -47215.922551188734
$ws.Range("E28").Value = 29836.299613162904
$ws.Range("F28").Value = 0.31262597079141563

$ws.Range("C29").Value = -36062.88023059594
$ws.Range("D29").Value = 70649.100601090759
$ws.Range("E29").Value = 177361.08143277746
$ws.Range("F29").Value = 0.25857075570236077

$ws.Range("C30").Value = 13458.108510023478
$ws.Range("D30").Value = 117865.02315227949
$ws.Range("E30").Value = 222271.93779453551
$ws.Range("F30").Value = 0.023150607221503416

# --- New: one-way ANOVA, dominant frequency (rows 32-36) ----------------
$ws.Range("A32").Value = "one-way ANOVA, dominant frequency"

$ws.Range("A33").Value = "Source"
$ws.Range("B33").Value = "SS"
$ws.Range("C33").Value = "df"
$ws.Range("D33").Value = "MS"
$ws.Range("E33").Value = "F"
$ws.Range("F33").Value = "Prob>F"

$ws.Range("A34").Value = "Groups"
$ws.Range("B34").Value = 181.91019926456897
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 90.955099632284487
$ws.Range("E34").Value = 2.3778586898253011
$ws.Range("F34").Value = 0.10030954804940632

$ws.Range("A35").Value = "Error"
$ws.Range("B35").Value = 2639.3081731398906
$ws.Range("C35").Value = 69
$ws.Range("D35").Value = 38.25084308898392

$ws.Range("A36").Value = "Total"
$ws.Range("B36").Value = 2821.2183724044598
$ws.Range("C36").Value = 71

# --- New: Multiple Comparison (Tukey-Kramer method), dominant frequency
#     (rows 38-42) --------------------------------------------------------
$ws.Range("A38").Value = "Multiple Comparison (Tukey-Kramer method), dominant frequency"

$ws.Range("A39").Value = "Group"
$ws.Range("B39").Value = "Group"
$ws.Range("F39").Value = "p-value"

$ws.Range("A40").Value = 1
$ws.Range("B40").Value = 2
$ws.Range("C40").Value = -0.41542396957261962
$ws.Range("D40").Value = 3.3909505208333357
$ws.Range("E40").Value = 7.1973250112392915
$ws.Range("F40").Value = 0.090451079898563491

$ws.Range("A41").Value = 1
$ws.Range("B41").Value = 3
$ws.Range("C41").Value = -2.5260506183858746
$ws.Range("D41").Value = 2.7455139160156286
$ws.Range("E41").Value = 8.0170784504171309
$ws.Range("F41").Value = 0.42963217002346643

$ws.Range("A42").Value = 2
$ws.Range("B42").Value = 3
$ws.Range("C42").Value = -5.8031310275508989
$ws.Range("D42").Value = -0.64543660481770715
$ws.Range("E42").Value = 4.5122578179154846
$ws.Range("F42").Value = 0.95171280623759369

# --- Mirror the author's final view state: scrolled/selected onto the new
#     Tukey-Kramer table at the bottom of the sheet ----------------------
$ws.Range("A40:F42").Select()
